# feat: add 2022-Q4 data
#
# 1. Insert a new worksheet "2022-Q4" right after "总计" (i.e. before the
#    existing "2022-Q2" sheet), populated with the quarterly fund-holding
#    detail table.
# 2. Insert a new row at the top of the "总计" (summary) sheet's data
#    block recording the 2022-Q4 totals, pushing every other quarter down
#    by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value that *looks* numeric (leading zeros, decimal
# text such as fund codes/percentages) into a cell while forcing Excel
# to keep it as literal text instead of silently parsing it into a
# number. Writing with a leading apostrophe forces text, and the
# subsequent ClearFormats() drops the quote-prefix style bit again
# (the stored type remains text) so we don't leave stray formatting
# behind on cells that should use the sheet's default style.
# ---------------------------------------------------------------------
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.ClearFormats()
}

# =======================================================================
# Part 1: new "2022-Q4" worksheet
# =======================================================================
# Duplicate the existing "2022-Q2" sheet (placed right before it) so the
# new sheet inherits identical sheet properties / page setup / header
# and column-A formatting, then overwrite its single data row with the
# four 2022-Q4 fund rows.
$refSheet = $wb.Worksheets.Item("2022-Q2")
$refSheet.Copy($refSheet)
$newSheet = $wb.Worksheets.Item("2022-Q2 (2)")
$newSheet.Name = "2022-Q4"

$newSheet.Range("A2:H2").ClearContents()
$newSheet.Range("A2").Copy()
$newSheet.Range("A2:A5").PasteSpecial(-4122)

$fundRows = @(
    @{ A = 0; B = "012073"; C = "华安均衡优选混合A";    D = "7.16"; E = "85.60"; F = "4.85"; G = "0.3473"; H = 2 },
    @{ A = 1; B = "001581"; C = "华安沪港深通精选混合A"; D = "5.03"; E = "83.24"; F = "5.32"; G = "0.2676"; H = 4 },
    @{ A = 2; B = "012074"; C = "华安均衡优选混合C";    D = "0.22"; E = "85.60"; F = "4.85"; G = "0.0107"; H = 2 },
    @{ A = 3; B = "016289"; C = "华安沪港深通精选混合C"; D = "0.09"; E = "83.24"; F = "5.32"; G = "0.0048"; H = 4 }
)

$r = 2
foreach ($row in $fundRows) {
    $newSheet.Range("A$r").Value = $row.A
    Set-TextValue $newSheet.Range("B$r") $row.B
    $newSheet.Range("C$r").Value = $row.C
    Set-TextValue $newSheet.Range("D$r") $row.D
    Set-TextValue $newSheet.Range("E$r") $row.E
    Set-TextValue $newSheet.Range("F$r") $row.F
    Set-TextValue $newSheet.Range("G$r") $row.G
    $newSheet.Range("H$r").Value = $row.H
    $r = $r + 1
}

# =======================================================================
# Part 2: update the "总计" summary sheet
# =======================================================================
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()

# Reuse the existing bold/centered/bordered style for the new A2 label
# cell (same style already used by the other rows in column A).
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.63

# Restore the originally-active tab ("2020-Q4", now shifted one slot to
# the right by the inserted sheet) so we don't leave the brand-new sheet
# marked as selected.
$wb.Worksheets.Item("2020-Q4").Activate()
